# Automatic map update (Actualización automática del mapa)
#
# The source data feed removed a handful of stale/resolved case rows from
# the "Optical_Power" sheet. Deleting full rows causes every row below the
# deleted one to shift up, which is exactly the effect observed between
# the original workbook (data through row 55) and the updated one (data
# through row 48).
#
# Rows removed (identified by their original row numbers / "Caso" values):
#   row 14 (Caso -51)   row 15 (Caso -212)  row 35 (Caso -416)
#   row 39 (Caso -437)  row 46 (Caso -451)  row 49 (Caso -502)
#   row 50 (Caso -506)
#
# Deleting from the bottom up keeps the remaining row numbers stable while
# we work, so the row list below is processed in descending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(50, 49, 46, 39, 35, 15, 14)

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
